$d = $word.ActiveDocument

# The document contains four paragraphs whose visible text is exactly
# "<id>pNNN</id>" split across three runs:
#   run 1: "<id>"   (Courier New, color 7f6000, size 9pt)  -- tagging colour
#   run 2: "pNNN"   (default body font, color 000000)
#   run 3: "</id>"  (same Courier New tagging colour as run 1)
# The edit merges the three runs into a single run (keeping run 1's
# formatting) whose text is the full "<id>pNNN</id>" string, while
# leaving the paragraph's trailing (empty) run untouched.

$ids = @("p031v_1", "p031v_2", "p031v_3", "p031v_4")

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $r = $p.Range

    foreach ($id in $ids) {
        $full = "<id>" + $id + "</id>"
        $paraText = $r.Text.TrimEnd([char]13, [char]7)
        if ($paraText -eq $full) {
            # Grab the formatting of the first run ("<id>") so the merged
            # run keeps its font/colour/size.
            $tagFont = $d.Range($r.Start, $r.Start + 4).Font

            $r.Find.Replacement.Text = $full
            $r.Find.Replacement.Font = $tagFont
            $r.Find.Execute($full, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $full, 2) | Out-Null
        }
    }
}
